# Aluminium construction, new pcb backplane with MPPC array
#
# - A2 label "plasticThickness" -> "caseThickness" (new parameter row)
# - A6 label stays "aerogelHolderDepth" (shared-string re-pointed after reorder)
# - B9 value 5.75 -> 4.7300000000000004
# - Active cell/selection moves from B10 to B9
# - Page setup: A4 portrait

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "plasticThickness" parameter row to "caseThickness"
$ws.Range("A2").Value = "caseThickness"

# Update the depth/backplane thickness value
$ws.Range("B9").Value = 4.7300000000000004

# Move the active selection to B9
$ws.Range("B9").Select()

# Configure page setup for printing (A4, portrait)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
